$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 160
$ws.Range("F6").Value = 314
$ws.Range("F7").Value = 5685
$ws.Range("F9").Value = 7665
$ws.Range("F10").Value = 396
$ws.Range("F11").Value = 68
$ws.Range("F12").Value = 53
$ws.Range("F13").Value = 3865
$ws.Range("F14").Value = 23
$ws.Range("F16").Value = 204
$ws.Range("F21").Value = 610
$ws.Range("F22").Value = 3897
$ws.Range("F23").Value = 135
$ws.Range("F25").Value = 5320
$ws.Range("F27").Value = 2110
$ws.Range("F28").Value = 132
$ws.Range("F29").Value = 354
$ws.Range("F30").Value = 7919
$ws.Range("F32").Value = 174
$ws.Range("F33").Value = 2201
$ws.Range("F34").Value = 2198
$ws.Range("F35").Value = 1335
$ws.Range("F36").Value = 1300
$ws.Range("F38").Value = 22
$ws.Range("F39").Value = 270
$ws.Range("F40").Value = 249
$ws.Range("F42").Value = 1178
$ws.Range("F44").Value = 35
$ws.Range("F45").Value = 1328
$ws.Range("F46").Value = 2089
$ws.Range("F47").Value = 130
$ws.Range("F48").Value = 222

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F20").Value = 13

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 573
$ws.Range("F3").Value = 752

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 160
$ws.Range("F5").Value = 573
$ws.Range("F6").Value = 752
$ws.Range("F7").Value = 314
$ws.Range("F8").Value = 5685
$ws.Range("F9").Value = 7665
$ws.Range("F10").Value = 397
$ws.Range("F11").Value = 3865
$ws.Range("F12").Value = 23
$ws.Range("F14").Value = 204
$ws.Range("F20").Value = 610
$ws.Range("F21").Value = 3897
$ws.Range("F23").Value = 135
$ws.Range("F25").Value = 5320
$ws.Range("F27").Value = 2110
$ws.Range("F28").Value = 132
$ws.Range("F29").Value = 354
$ws.Range("F30").Value = 7919
$ws.Range("F32").Value = 174
$ws.Range("F33").Value = 2201
$ws.Range("F34").Value = 2198
$ws.Range("F35").Value = 1335
$ws.Range("F36").Value = 1300
$ws.Range("F37").Value = 270
$ws.Range("F38").Value = 249
$ws.Range("F40").Value = 1178
$ws.Range("F42").Value = 35
$ws.Range("F43").Value = 1328
$ws.Range("F44").Value = 2089
$ws.Range("F45").Value = 130
$ws.Range("F47").Value = 222
$ws.Range("F48").Value = 13
